$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 672, pushing the existing rows 672:728 down to 673:729
$ws.Rows.Item(672).Insert()

# Populate the newly inserted row 672 with the new data entry
$ws.Cells.Item(672, 1).Value = 10
$ws.Cells.Item(672, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(672, 3).Value = "La Araucanía"
$ws.Cells.Item(672, 4).Value = 45223
$ws.Cells.Item(672, 5).Value = 9
$ws.Cells.Item(672, 6).Value = 100112037
$ws.Cells.Item(672, 7).Value = "Cebollín"
$ws.Cells.Item(672, 8).Value = "Sin especificar"
$ws.Cells.Item(672, 9).Value = "Primera"
$ws.Cells.Item(672, 10).Value = 30
$ws.Cells.Item(672, 11).Value = 8000
$ws.Cells.Item(672, 12).Value = 8000
$ws.Cells.Item(672, 13).Value = 8000
$ws.Cells.Item(672, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(672, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(672, 16).Value = 667
$ws.Cells.Item(672, 17).Value = 12
$ws.Cells.Item(672, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by the rest of column D
$ws.Cells.Item(672, 4).NumberFormat = $ws.Cells.Item(673, 4).NumberFormat
